$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = "4VDyVobDnbtrJVlHnSdP"
$ws.Range("D6").Value = "R8RgpUn4Uv52BwTDBxTm"

$ws.Range("C12").Value = "B2YjkjRDy4C3LfWOTAVE"
$ws.Range("D12").Value = "GjRFa1PPsr1KtgAY57cz"

$ws.Range("C18").Value = "tIeFMO0HT1YdkchZn5lA"
$ws.Range("D18").Value = "xpABrtBO6e5WF7JalJ6I"

$ws.Range("E24").Value = "TNl9SZJTalf6qjJp4S8Q"
$ws.Range("F24").Value = "6ClGRU9A57BvdmRQPhA7"

$ws.Range("F24").Select()
$excel.ActiveWindow.ScrollRow = 4
